# "add read, update and delete for zones"
#
# The "Test Results" sheet tracks pass/fail (TRUE/FALSE) of CRUD tests for
# Device/Zone/Category records. Rows 2-10 are the Zone (Z01-Z09) rows, and
# column B is "Create Test Passed". With read/update/delete now implemented
# for zones, the Create test for each zone row now passes as well, so flip
# B2:B10 from FALSE to TRUE. The active selection on the sheet is also left
# on F4, reflecting where the author's cursor ended up while making the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Create Test Passed -> TRUE for every Zone row (Z01..Z09)
$ws.Range("B2:B10").Value = $true

# Leave the sheet active with the selection where the author left it
$ws.Activate()
$ws.Range("F4").Select()
